$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Enter the week 1 progress note for Anushka Chincholkar in column B, row 4
$ws.Range("B4").Value = "Studied basic concepts of Flutter toolkit"

# Reflect the active selection left after typing the entry
$ws.Range("B4").Select()
